$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 879, pushing existing rows 879-920 down to 880-921.
$ws.Rows(879).Insert()

# Populate the newly inserted row 879 with the new entry: 2026/02/26, Thursday, 17:00, rank 201.
# Format column A as text first so the date string is kept literally (not converted to a date serial),
# matching the existing cells in column A which are plain text.
$ws.Range("A879").NumberFormat = "@"
$ws.Range("A879").Value = "2026/02/26"
$ws.Range("A879").Style = "Normal"

$ws.Range("B879").Value = "木"

$ws.Range("C879").Value = 17

$ws.Range("D879").Value = 201
